$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.041.90"
$ws.Range("D3").Value = "1.791.88"
$ws.Range("E3").Value = "  -0.70%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").Value = "'223.44"
$ws.Range("E5").Value = "  -0.53%  "
$ws.Range("E6").Value = "  -0.63%  "
$ws.Range("E7").Value = "  +0.10%  "
$ws.Range("D8").Value = "'32.33"
$ws.Range("E8").Value = "  -0.93%  "
$ws.Range("D9").Value = "'0.284"
$ws.Range("E9").Value = "  -2.38%  "
$ws.Range("D10").Value = "'0.0707"
$ws.Range("E10").Value = "  -0.54%  "
$ws.Range("E11").Value = "  +0.02%  "
$ws.Range("D12").Value = "2.052.11"
$ws.Range("E12").Value = "  -0.63%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.791.91"
$ws.Range("E13").Value = "  -0.95%  "
$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").Value = "'10.89"
$ws.Range("E14").Value = "  -2.02%  "
$ws.Range("D15").Value = "'0.624"
$ws.Range("E15").Value = "  -3.01%  "
$ws.Range("D16").Value = "34.082.95"
$ws.Range("D17").Value = "'4.15"
$ws.Range("E17").Value = "  -4.37%  "
$ws.Range("D18").Value = "'67.92"
$ws.Range("D19").Value = "'243.10"
$ws.Range("E19").Value = "  -4.18%  "
$ws.Range("D20").Value = "0.0₃0781"
$ws.Range("E20").Value = "  -2.55%  "
$ws.Range("E21").Value = "  +0.10%  "
$ws.Range("D22").Value = "'10.68"
$ws.Range("E22").Value = "  -4.10%  "
$ws.Range("D23").Value = "'4.07"
$ws.Range("E23").Value = "  -4.46%  "
$ws.Range("E24").Value = "  -2.95%  "
$ws.Range("D25").Value = "'158.85"
$ws.Range("E25").Value = "  -1.72%  "
$ws.Range("D26").Value = "'16.25"
$ws.Range("E26").Value = "  -1.31%  "
$ws.Range("D27").Value = "'7.00"
$ws.Range("E27").Value = "  -2.26%  "
$ws.Range("E28").Value = "  -2.18%  "
$ws.Range("E29").Value = "  +0.06%  "
$ws.Range("D30").Value = "'0.0518"
$ws.Range("E30").Value = "  -2.26%  "
$ws.Range("E31").Value = "  +0.31%  "
$ws.Range("D32").Value = "'3.66"
$ws.Range("E32").Value = "  -3.84%  "
$ws.Range("D33").Value = "'3.49"
$ws.Range("E33").Value = "  -4.14%  "
$ws.Range("D34").Value = "'1.81"
$ws.Range("E34").Value = "  -4.61%  "
$ws.Range("D35").Value = "1.384.09"
$ws.Range("E35").Value = "  -3.62%  "
$ws.Range("D36").Value = "'0.646"
$ws.Range("E36").Value = "  +0.05%  "
$ws.Range("E37").Value = "  -2.01%  "
$ws.Range("D38").Value = "'0.0184"
$ws.Range("E38").Value = "  -4.17%  "
$ws.Range("D39").Value = "'79.42"
$ws.Range("E39").Value = "  -6.57%  "
$ws.Range("E40").Value = "  +0.22%  "
$ws.Range("B41").Value = "MXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D41").Value = "'2.70"
$ws.Range("E41").Value = "  -3.59%  "
$ws.Range("B42").Value = "ARBITRUM"
$ws.Range("C42").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D42").Value = "'0.913"
$ws.Range("E42").Value = "  -4.56%  "
$ws.Range("D43").Value = "'2.16"
$ws.Range("E43").Value = "  -0.25%  "
$ws.Range("D44").Value = "0.0₆0138"
$ws.Range("E44").Value = "  +8.47%  "
$ws.Range("E45").Value = "  +0.39%  "
$ws.Range("B46").Value = "WEMIXToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D46").Value = "'1.05"
$ws.Range("E46").Value = "  -0.63%  "
$ws.Range("B47").Value = "Quant"
$ws.Range("C47").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D47").Value = "'107.26"
$ws.Range("E47").Value = "  +0.86%  "
$ws.Range("D48").Value = "'5.85"
$ws.Range("E48").Value = "  -3.37%  "
$ws.Range("D49").Value = "1.951.91"
$ws.Range("E49").Value = "  -0.33%  "
$ws.Range("B50").Value = "PaxDollar"
$ws.Range("C50").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D50").Value = "'1.00"
$ws.Range("E50").Value = "  +0.06%  "
$ws.Range("B51").Value = "InjectiveProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D51").Value = "'11.95"
$ws.Range("E51").Value = "  -2.68%  "
